$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) price cells to remain text, matching the original inlineStr type,
# so values like "73.733.16" or "0.999" are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '73.733.16'
$ws.Cells.Item(2, 5).Value = '  +7.27%  '
$ws.Cells.Item(3, 4).Value = '2.616.53'
$ws.Cells.Item(3, 5).Value = '  +7.27%  '
$ws.Cells.Item(4, 5).Value = '  +0.02%  '
$ws.Cells.Item(5, 4).Value = '185.98'
$ws.Cells.Item(5, 5).Value = '  +13.88%  '
$ws.Cells.Item(6, 4).Value = '581.18'
$ws.Cells.Item(6, 5).Value = '  +3.73%  '
$ws.Cells.Item(7, 5).Value = '  -0.06%  '
$ws.Cells.Item(8, 4).Value = '0.529'
$ws.Cells.Item(8, 5).Value = '  +4.44%  '
$ws.Cells.Item(9, 4).Value = '0.197'
$ws.Cells.Item(9, 5).Value = '  +16.23%  '
$ws.Cells.Item(10, 4).Value = '2.613.23'
$ws.Cells.Item(10, 5).Value = '  +7.23%  '
$ws.Cells.Item(11, 4).Value = '0.163'
$ws.Cells.Item(11, 5).Value = '  +1.14%  '
$ws.Cells.Item(12, 4).Value = '0.356'
$ws.Cells.Item(12, 5).Value = '  +7.48%  '
$ws.Cells.Item(13, 4).Value = '4.68'
$ws.Cells.Item(13, 5).Value = '  +1.55%  '
$ws.Cells.Item(14, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(14, 4).Value = '3.120.35'
$ws.Cells.Item(14, 5).Value = '  +8.09%  '
$ws.Cells.Item(15, 4).Value = '73.748.61'
$ws.Cells.Item(15, 5).Value = '  +7.47%  '
$ws.Cells.Item(16, 2).Value = 'ShibaInu'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(16, 4).Value = '0.0000188'
$ws.Cells.Item(16, 5).Value = '  +5.05%  '
$ws.Cells.Item(17, 4).Value = '26.25'
$ws.Cells.Item(17, 5).Value = '  +12.51%  '
$ws.Cells.Item(18, 4).Value = '2.620.15'
$ws.Cells.Item(18, 5).Value = '  +7.28%  '
$ws.Cells.Item(19, 4).Value = '9.08'
$ws.Cells.Item(19, 5).Value = '  +29.78%  '
$ws.Cells.Item(20, 4).Value = '11.80'
$ws.Cells.Item(20, 5).Value = '  +11.37%  '
$ws.Cells.Item(21, 4).Value = '365.11'
$ws.Cells.Item(21, 5).Value = '  +7.78%  '
$ws.Cells.Item(22, 4).Value = '2.27'
$ws.Cells.Item(22, 5).Value = '  +16.72%  '
$ws.Cells.Item(23, 4).Value = '4.05'
$ws.Cells.Item(23, 5).Value = '  +5.47%  '
$ws.Cells.Item(24, 4).Value = '0.999'
$ws.Cells.Item(24, 5).Value = '  -0.12%  '
$ws.Cells.Item(25, 4).Value = '69.43'
$ws.Cells.Item(25, 5).Value = '  +5.72%  '
$ws.Cells.Item(26, 4).Value = '4.11'
$ws.Cells.Item(26, 5).Value = '  +9.19%  '
$ws.Cells.Item(27, 4).Value = '9.25'
$ws.Cells.Item(27, 5).Value = '  +10.98%  '
$ws.Cells.Item(28, 4).Value = '2.762.39'
$ws.Cells.Item(28, 5).Value = '  +7.58%  '
$ws.Cells.Item(29, 5).Value = '  +0.53%  '
$ws.Cells.Item(30, 4).Value = '0.0₃0936'
$ws.Cells.Item(30, 5).Value = '  +13.84%  '
$ws.Cells.Item(31, 4).Value = '522.17'
$ws.Cells.Item(31, 5).Value = '  +20.75%  '
$ws.Cells.Item(32, 4).Value = '1.37'
$ws.Cells.Item(32, 5).Value = '  +15.55%  '
$ws.Cells.Item(33, 4).Value = '7.65'
$ws.Cells.Item(33, 5).Value = '  +6.61%  '
$ws.Cells.Item(34, 4).Value = '1.74'
$ws.Cells.Item(34, 5).Value = '  +9.02%  '
$ws.Cells.Item(35, 4).Value = '0.999'
$ws.Cells.Item(35, 5).Value = '  -0.03%  '
$ws.Cells.Item(36, 4).Value = '161.35'
$ws.Cells.Item(36, 5).Value = '  +1.50%  '
$ws.Cells.Item(37, 5).Value = '  +9.53%  '
$ws.Cells.Item(38, 4).Value = '19.05'
$ws.Cells.Item(38, 5).Value = '  +5.96%  '
$ws.Cells.Item(39, 4).Value = '19.26'
$ws.Cells.Item(39, 5).Value = '  +1.34%  '
$ws.Cells.Item(40, 5).Value = '  +0.07%  '
$ws.Cells.Item(41, 4).Value = '4.89'
$ws.Cells.Item(41, 5).Value = '  +11.98%  '
$ws.Cells.Item(42, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(42, 4).Value = '0.325'
$ws.Cells.Item(42, 5).Value = '  +8.06%  '
$ws.Cells.Item(43, 2).Value = 'Stacks'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(43, 4).Value = '1.65'
$ws.Cells.Item(43, 5).Value = '  +9.48%  '
$ws.Cells.Item(44, 4).Value = '161.54'
$ws.Cells.Item(44, 5).Value = '  +24.16%  '
$ws.Cells.Item(45, 4).Value = '2.36'
$ws.Cells.Item(45, 5).Value = '  +13.63%  '
$ws.Cells.Item(46, 4).Value = '1.17'
$ws.Cells.Item(46, 5).Value = '  +9.06%  '
$ws.Cells.Item(47, 4).Value = '38.86'
$ws.Cells.Item(47, 5).Value = '  +3.61%  '
$ws.Cells.Item(48, 4).Value = '0.0845'
$ws.Cells.Item(48, 5).Value = '  +17.34%  '
$ws.Cells.Item(49, 4).Value = '3.59'
$ws.Cells.Item(49, 5).Value = '  +8.05%  '
$ws.Cells.Item(50, 4).Value = '0.523'
$ws.Cells.Item(50, 5).Value = '  +8.38%  '
$ws.Cells.Item(51, 4).Value = '20.74'
$ws.Cells.Item(51, 5).Value = '  +22.87%  '

# Restore default style on column D so only values changed (no residual number format).
$ws.Range("D2:D51").Style = "Normal"
